$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (label unchanged)
$ws.Range("B2").Value = 0.205383640252233
$ws.Range("C2").Value = 0.2053836402522327
$ws.Range("D2").Value = 0.205383640252233

# Row 3 - RandomForestRegressor (label unchanged)
$ws.Range("B3").Value = 0.02398204792323529
$ws.Range("C3").Value = 0.02477136214647617
$ws.Range("D3").Value = 0.04945056209067746

# Row 4 - label changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02256747169237085
$ws.Range("C4").Value = 0.0251526531094253
$ws.Range("D4").Value = 0.05651174749734465

# Row 5 - label changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.01985224583206885
$ws.Range("C5").Value = 0.02124911211142868
$ws.Range("D5").Value = 0.02231792118752034
